# Update report header: issue number and week-covering dates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "Volume 31   Number  28"
$ws.Range("C9").Value = "Report Covering the Week  7/8/2024  Through  7/14/2024"

# Cells that were previously placeholder text ("0" / "***.*") and must become
# real numeric values need their number format restored so they store as
# numbers instead of text.
$ws.Range("C14").NumberFormat = "#,##0"
$ws.Range("D14").NumberFormat = "#,##0"
$ws.Range("E14").NumberFormat = "#,##0.0;""-""#,##0.0"

# New weekly crime-complaint figures (rows 14-31) and the Traffic Fatalities
# row (row 33). Columns:
#   C/D = Week to Date (2024/2023), E = % Chg
#   F/G = 28 Day (2024/2023),       H = % Chg
#   I/J = Year to Date (2024/2023), K = % Chg
#   L   = 2 Year % Chg
#   M   = 14 Year % Chg
#   N   = 31 Year % Chg
$rows = @(
  @{ Row=14; C=2; D=1; E=100; F=4; G=4; H=0; I=39; J=36; K=8.333333333333; L=-9.302325581395; M=-47.297297297297; N=-83.817427385892 },
  @{ Row=15; C=2; D=7; E=-71.428571428571; F=19; G=24; H=-20.833333333333; I=135; J=127; K=6.299212598425; L=10.655737704918; M=22.727272727272; N=-60.176991150442 },
  @{ Row=16; C=43; D=42; E=2.380952380952; F=166; G=175; H=-5.142857142857; I=1278; J=1262; K=1.267828843106; L=-6.647187728268; M=-28.483491885842; N=-85.06660434681 },
  @{ Row=17; C=83; D=97; E=-14.432989690721; F=367; G=389; H=-5.655526992287; I=2315; J=2294; K=0.915431560592; L=6.338998621956; M=30.717108977978; N=-48.302813756141 },
  @{ Row=18; C=26; D=41; E=-36.585365853658; F=118; G=161; H=-26.708074534161; I=992; J=1085; K=-8.571428571428; L=-22.317932654659; M=-35.542560103963; N=-84.193753983428 },
  @{ Row=19; C=96; D=138; E=-30.434782608695; F=410; G=496; H=-17.338709677419; I=2784; J=3101; K=-10.222508868107; L=-8.811005568293; M=33.205741626794; N=-18.881118881118 },
  @{ Row=20; C=35; D=50; E=-30; F=143; G=147; H=-2.721088435374; I=885; J=928; K=-4.633620689655; L=-4.427645788336; M=16.908850726552; N=-82.377538829151 },
  @{ Row=21; C=287; D=376; E=-23.670212765957; F=1227; G=1396; H=-12.106017191977; I=8428; J=8833; K=-4.585078682214; L=-6.010928961748; M=3.690944881889; N=-70.267409863825 },
  @{ Row=22; C=7; D=7; E=0; F=26; G=27; H=-3.703703703703; I=150; J=160; K=-6.25; L=-25.742574257425; M=-31.818181818181; N="***.*" },
  @{ Row=23; C=34; D=34; E=0; F=104; G=130; H=-20; I=797; J=879; K=-9.328782707622; L=-0.870646766169; M=35.544217687074; N="***.*" },
  @{ Row=24; C=251; D=265; E=-5.283018867924; F=999; G=1052; H=-5.038022813688; I=6439; J=6600; K=-2.439393939393; L=-6.923966464296; M=21.927665214921; N="***.*" },
  @{ Row=25; C=129; D=115; E=12.173913043478; F=469; G=411; H=14.111922141119; I=2849; J=2642; K=7.83497350492; L=0.849557522123; M="***.*"; N="***.*" },
  @{ Row=26; C=125; D=137; E=-8.759124087591; F=536; G=555; H=-3.423423423423; I=3418; J=3271; K=4.49403852033; L=4.175556232855; M=-19.859320046893; N="***.*" },
  @{ Row=27; C=3; D=9; E=-66.666666666666; F=22; G=37; H=-40.54054054054; I=195; J=192; K=1.5625; L=1.036269430051; M="***.*"; N="***.*" },
  @{ Row=28; C=18; D=12; E=50; F=65; G=43; H=51.162790697674; I=345; J=322; K=7.142857142857; L=4.545454545454; M="***.*"; N="***.*" },
  @{ Row=29; C=12; D=4; E=200; F=26; G=19; H=36.842105263157; I=132; J=123; K=7.317073170731; L=-26.666666666666; M=-50.375939849624; N=-87.058823529411 },
  @{ Row=30; C=10; D=4; E=150; F=23; G=18; H=27.777777777777; I=115; J=108; K=6.481481481481; L=-24.342105263157; M=-44.976076555023; N=-87.486398258977 },
  @{ Row=31; C="0"; D=2; E=-100; F=2; G=9; H=-77.777777777777; I=40; J=42; K=-4.761904761904; L=-2.439024390243; M="***.*"; N="***.*" },
  @{ Row=33; C="0"; D="0"; E="***.*"; F="0"; G=1; H=-100; I=13; J=9; K=44.444444444444; L=-13.333333333333; M="***.*"; N="***.*" }
)

foreach ($r in $rows) {
    $rowNum = $r.Row
    foreach ($col in @("C","D","E","F","G","H","I","J","K","L","M","N")) {
        if ($r.ContainsKey($col)) {
            $addr = "$col$rowNum"
            $ws.Range($addr).Value = $r[$col]
        }
    }
}
